$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: title + generated-on date -----------------------------------
$ws.Range("A1").Value = "PRODUCTION STATUS REPORT"
$ws.Range("B2").Value = "21/1/2019"

# Helper: write a value that must be stored as TEXT even when it looks like
# a pure number (quantities / done-counts are text cells in this report),
# by using Excel's quote-prefix ('-) entry so it is not reinterpreted as a
# numeric literal.
function Set-TextValue($rng, $val) {
    if ($val -match '^-?\d+(\.\d+)?$') {
        $rng.Value = "'" + $val
    } else {
        $rng.Value = $val
    }
}

# --- Data rows 6-9 (job-slip detail) --------------------------------------
# Row 6: Dream uniforms / ARKAN SECURITY UNFIORM - Shirt
$ws.Range("A6").Value = 1
Set-TextValue $ws.Range("B6") "Dream uniforms"
Set-TextValue $ws.Range("C6") "lpo/Dream uniforms/73315"
Set-TextValue $ws.Range("D6") "5857-5956"
Set-TextValue $ws.Range("E6") "Shirt"
Set-TextValue $ws.Range("F6") "ARKAN SECURITY UNFIORM"
Set-TextValue $ws.Range("G6") "100"
Set-TextValue $ws.Range("H6") "0"
Set-TextValue $ws.Range("I6") " cut on "
Set-TextValue $ws.Range("J6") "0"
Set-TextValue $ws.Range("K6") " stitched on "
Set-TextValue $ws.Range("L6") "0"
Set-TextValue $ws.Range("M6") " finalized on "
Set-TextValue $ws.Range("N6") "100"
Set-TextValue $ws.Range("O6") "riaz packing on 21/1/2019--20:14"
Set-TextValue $ws.Range("P6") "0"

# Row 7: Dream uniforms / ARKAN SECURITY UNFIORM - Trouser
$ws.Range("A7").Value = 2
Set-TextValue $ws.Range("B7") "Dream uniforms"
Set-TextValue $ws.Range("C7") "lpo/Dream uniforms/73315"
Set-TextValue $ws.Range("D7") "5957-6006"
Set-TextValue $ws.Range("E7") "Trouser"
Set-TextValue $ws.Range("F7") "ARKAN SECURITY UNFIORM"
Set-TextValue $ws.Range("G7") "50"
Set-TextValue $ws.Range("H7") "0"
Set-TextValue $ws.Range("I7") " cut on "
Set-TextValue $ws.Range("J7") "0"
Set-TextValue $ws.Range("K7") " stitched on "
Set-TextValue $ws.Range("L7") "0"
Set-TextValue $ws.Range("M7") " finalized on "
Set-TextValue $ws.Range("N7") "50"
Set-TextValue $ws.Range("O7") "mujtaba packing on 21/1/2019--20:12"
Set-TextValue $ws.Range("P7") "0"

# Row 8: Dream uniforms / ARKAN SECURITY UNFIORM - Shirt
$ws.Range("A8").Value = 3
Set-TextValue $ws.Range("B8") "Dream uniforms"
Set-TextValue $ws.Range("C8") "lpo/Dream uniforms/73315"
Set-TextValue $ws.Range("D8") "6007-6106"
Set-TextValue $ws.Range("E8") "Shirt"
Set-TextValue $ws.Range("F8") "ARKAN SECURITY UNFIORM"
Set-TextValue $ws.Range("G8") "100"
Set-TextValue $ws.Range("H8") "0"
Set-TextValue $ws.Range("I8") " cut on "
Set-TextValue $ws.Range("J8") "0"
Set-TextValue $ws.Range("K8") " stitched on "
Set-TextValue $ws.Range("L8") "0"
Set-TextValue $ws.Range("M8") " finalized on "
Set-TextValue $ws.Range("N8") "24"
Set-TextValue $ws.Range("O8") "ali packing on 21/1/2019--20:24"
Set-TextValue $ws.Range("P8") "76"

# Row 9: Dream uniforms / ARKAN SECURITY UNFIORM - Trouser
$ws.Range("A9").Value = 4
Set-TextValue $ws.Range("B9") "Dream uniforms"
Set-TextValue $ws.Range("C9") "lpo/Dream uniforms/73315"
Set-TextValue $ws.Range("D9") "6107-6156"
Set-TextValue $ws.Range("E9") "Trouser"
Set-TextValue $ws.Range("F9") "ARKAN SECURITY UNFIORM"
Set-TextValue $ws.Range("G9") "50"
Set-TextValue $ws.Range("H9") "0"
Set-TextValue $ws.Range("I9") " cut on "
Set-TextValue $ws.Range("J9") "0"
Set-TextValue $ws.Range("K9") " stitched on "
Set-TextValue $ws.Range("L9") "0"
Set-TextValue $ws.Range("M9") " finalized on "
Set-TextValue $ws.Range("N9") "0"
Set-TextValue $ws.Range("O9") " packing on "
Set-TextValue $ws.Range("P9") "50"

# --- Row 10 is dropped entirely; this also shifts the old totals row 12 --
# up to become row 11, and fixes the sheet dimension automatically.
$ws.Rows(10).Delete()

# --- Totals row (was row 12, now row 11 after the delete/shift) ----------
Set-TextValue $ws.Range("G11") "Total Cut"
$ws.Range("H11").Value = 0
Set-TextValue $ws.Range("I11") "Total Stitched"
$ws.Range("J11").Value = 0
Set-TextValue $ws.Range("K11") "Total Finished"
$ws.Range("L11").Value = 0
Set-TextValue $ws.Range("M11") "Total Packed"
$ws.Range("N11").Value = 174
Set-TextValue $ws.Range("O11") "Total Delivered"
$ws.Range("P11").Value = 174
